# Generate Report for Handoff
# Records a fresh "Latest Handoff Date(time)" timestamp for the most
# recently handed-off file (b6c45c8e-3ec2-4425-8c5d-8dfa3008d792) across
# the Overview roll-up sheet and each per-locale detail sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D7").Value = "2016-27-11 22:27:09"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E7").Value = "2016-03-11 22:27:06"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E7").Value = "2016-03-11 22:27:09"
